$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Cells.Item(74, 8).Value = 5225.375  # H74
$ws.Cells.Item(74, 9).Value = 4880.6  # I74
$ws.Cells.Item(74, 10).Value = 5800  # J74
$ws.Cells.Item(74, 11).Value = 4880.6  # K74
$ws.Cells.Item(74, 12).Value = 5800  # L74
$ws.Cells.Item(74, 13).Value = -3944.6  # M74
$ws.Cells.Item(74, 14).Value = -7672  # N74
# Row 76
$ws.Cells.Item(76, 8).Value = 4324.6875  # H76
$ws.Cells.Item(76, 10).Value = 4865.8335  # J76
$ws.Cells.Item(76, 12).Value = 4865.8335  # L76
$ws.Cells.Item(76, 14).Value = -5495.8335  # N76
# Row 77
$ws.Cells.Item(77, 8).Value = 5225.375  # H77
$ws.Cells.Item(77, 9).Value = 4880.6  # I77
$ws.Cells.Item(77, 10).Value = 5800  # J77
$ws.Cells.Item(77, 11).Value = 24403  # K77
$ws.Cells.Item(77, 12).Value = 29000  # L77
$ws.Cells.Item(77, 13).Value = -19723  # M77
$ws.Cells.Item(77, 14).Value = -38360  # N77
# Row 79
$ws.Cells.Item(79, 8).Value = 4324.6875  # H79
$ws.Cells.Item(79, 10).Value = 4865.8335  # J79
$ws.Cells.Item(79, 12).Value = 4865.8335  # L79
$ws.Cells.Item(79, 14).Value = -7049.8335  # N79
# Row 107
$ws.Cells.Item(107, 8).Value = 652.75  # H107
$ws.Cells.Item(107, 9).Value = 592.7  # I107
$ws.Cells.Item(107, 10).Value = 953  # J107
$ws.Cells.Item(107, 11).Value = 592.7  # K107
$ws.Cells.Item(107, 12).Value = 953  # L107
$ws.Cells.Item(107, 13).Value = 1327.3  # M107
$ws.Cells.Item(107, 14).Value = -4793  # N107

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 30184.293  # H32
$ws.Cells.Item(32, 9).Value = 36071.902  # I32
$ws.Cells.Item(32, 11).Value = 36071.902  # K32
$ws.Cells.Item(32, 13).Value = -35784.902  # M32
# Row 63
$ws.Cells.Item(63, 8).Value = 3098.8  # H63
$ws.Cells.Item(63, 9).Value = 2498  # I63
$ws.Cells.Item(63, 11).Value = 2498  # K63
$ws.Cells.Item(63, 13).Value = -1812  # M63
# Row 66
$ws.Cells.Item(66, 8).Value = 3098.8  # H66
$ws.Cells.Item(66, 9).Value = 2498  # I66
$ws.Cells.Item(66, 11).Value = 12490  # K66
$ws.Cells.Item(66, 13).Value = -9058  # M66
# Row 132
$ws.Cells.Item(132, 8).Value = 2513.625  # H132
$ws.Cells.Item(132, 9).Value = 2094.682  # I132
$ws.Cells.Item(132, 10).Value = 2868.1155  # J132
$ws.Cells.Item(132, 11).Value = 6284.045999999999  # K132
$ws.Cells.Item(132, 12).Value = 8604.3465  # L132
$ws.Cells.Item(132, 13).Value = -3754.045999999999  # M132
$ws.Cells.Item(132, 14).Value = -13664.3465  # N132

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Cells.Item(22, 8).Value = 426.8  # H22
$ws.Cells.Item(22, 9).Value = 426.8  # I22
$ws.Cells.Item(22, 10).Value = 0  # J22
$ws.Cells.Item(22, 11).Value = 426.8  # K22
$ws.Cells.Item(22, 12).Value = 0  # L22
$ws.Cells.Item(22, 13).Value = -253.8  # M22
$ws.Cells.Item(22, 14).ClearContents()  # N22
# Row 64
$ws.Cells.Item(64, 8).Value = 624.5  # H64
$ws.Cells.Item(64, 10).Value = 749  # J64
$ws.Cells.Item(64, 12).Value = 749  # L64
$ws.Cells.Item(64, 14).Value = -1199  # N64
# Row 67
$ws.Cells.Item(67, 8).Value = 624.5  # H67
$ws.Cells.Item(67, 10).Value = 749  # J67
$ws.Cells.Item(67, 12).Value = 749  # L67
$ws.Cells.Item(67, 14).Value = -2309  # N67
# Row 109
$ws.Cells.Item(109, 8).Value = 35000  # H109
$ws.Cells.Item(109, 10).Value = 35000  # J109
$ws.Cells.Item(109, 12).Value = 35000  # L109
$ws.Cells.Item(109, 14).Value = -37774  # N109
# Row 134
$ws.Cells.Item(134, 8).Value = 2406.2856  # H134
$ws.Cells.Item(134, 9).Value = 2224  # I134
$ws.Cells.Item(134, 11).Value = 6672  # K134
$ws.Cells.Item(134, 13).Value = -4137  # M134

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 3253.6667  # H31
$ws.Cells.Item(31, 9).Value = 2789.3333  # I31
$ws.Cells.Item(31, 11).Value = 2789.3333  # K31
$ws.Cells.Item(31, 13).Value = -2494.3333  # M31
# Row 34
$ws.Cells.Item(34, 8).Value = 3253.6667  # H34
$ws.Cells.Item(34, 9).Value = 2789.3333  # I34
$ws.Cells.Item(34, 11).Value = 2789.3333  # K34
$ws.Cells.Item(34, 13).Value = -2587.3333  # M34
# Row 132
$ws.Cells.Item(132, 8).Value = 2784.5518  # H132
$ws.Cells.Item(132, 9).Value = 2446.9285  # I132
$ws.Cells.Item(132, 11).Value = 7340.7855  # K132
$ws.Cells.Item(132, 13).Value = -4810.7855  # M132
# Row 134
$ws.Cells.Item(134, 8).Value = 4472.1665  # H134
$ws.Cells.Item(134, 9).Value = 3866.75  # I134
$ws.Cells.Item(134, 11).Value = 11600.25  # K134
$ws.Cells.Item(134, 13).Value = -9065.25  # M134

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 136
$ws.Cells.Item(136, 8).Value = 3063.4614  # H136
$ws.Cells.Item(136, 9).Value = 1176.6666  # I136
$ws.Cells.Item(136, 10).Value = 3309.5652  # J136
$ws.Cells.Item(136, 11).Value = 3529.9998  # K136
$ws.Cells.Item(136, 12).Value = 9928.695599999999  # L136
$ws.Cells.Item(136, 13).Value = 1570.0002  # M136
$ws.Cells.Item(136, 14).Value = -20128.6956  # N136

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 5512.0967  # H70
$ws.Cells.Item(70, 9).Value = 5175.3687  # I70
$ws.Cells.Item(70, 10).Value = 5660.884  # J70
$ws.Cells.Item(70, 11).Value = 5175.3687  # K70
$ws.Cells.Item(70, 12).Value = 5660.884  # L70
$ws.Cells.Item(70, 13).Value = -4905.3687  # M70
$ws.Cells.Item(70, 14).Value = -6200.884  # N70
# Row 73
$ws.Cells.Item(73, 8).Value = 5512.0967  # H73
$ws.Cells.Item(73, 9).Value = 5175.3687  # I73
$ws.Cells.Item(73, 10).Value = 5660.884  # J73
$ws.Cells.Item(73, 11).Value = 5175.3687  # K73
$ws.Cells.Item(73, 12).Value = 5660.884  # L73
$ws.Cells.Item(73, 13).Value = -4239.3687  # M73
$ws.Cells.Item(73, 14).Value = -7532.884  # N73
# Row 132
$ws.Cells.Item(132, 8).Value = 2738.4348  # H132
$ws.Cells.Item(132, 9).Value = 2691.125  # I132
$ws.Cells.Item(132, 11).Value = 8073.375  # K132
$ws.Cells.Item(132, 13).Value = -5543.375  # M132
# Row 133
$ws.Cells.Item(133, 8).Value = 0  # H133
$ws.Cells.Item(133, 10).Value = 0  # J133
$ws.Cells.Item(133, 12).Value = 0  # L133
$ws.Cells.Item(133, 14).ClearContents()  # N133

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 5092.12  # H7
$ws.Cells.Item(7, 9).Value = 4331.737  # I7
$ws.Cells.Item(7, 10).Value = 7500  # J7
$ws.Cells.Item(7, 11).Value = 4331.737  # K7
$ws.Cells.Item(7, 12).Value = 7500  # L7
$ws.Cells.Item(7, 13).Value = -4219.737  # M7
$ws.Cells.Item(7, 14).Value = -7724  # N7
# Row 68
$ws.Cells.Item(68, 8).Value = 3142.3333  # H68
$ws.Cells.Item(68, 9).Value = 3205.3635  # I68
$ws.Cells.Item(68, 10).Value = 3043.2856  # J68
$ws.Cells.Item(68, 11).Value = 3205.3635  # K68
$ws.Cells.Item(68, 12).Value = 3043.2856  # L68
$ws.Cells.Item(68, 13).Value = -2456.3635  # M68
$ws.Cells.Item(68, 14).Value = -4541.2856  # N68
# Row 71
$ws.Cells.Item(71, 8).Value = 3142.3333  # H71
$ws.Cells.Item(71, 9).Value = 3205.3635  # I71
$ws.Cells.Item(71, 10).Value = 3043.2856  # J71
$ws.Cells.Item(71, 11).Value = 16026.8175  # K71
$ws.Cells.Item(71, 12).Value = 15216.428  # L71
$ws.Cells.Item(71, 13).Value = -12282.8175  # M71
$ws.Cells.Item(71, 14).Value = -22704.428  # N71
# Row 82
$ws.Cells.Item(82, 8).Value = 0  # H82
$ws.Cells.Item(82, 9).Value = 0  # I82
$ws.Cells.Item(82, 10).Value = 0  # J82
$ws.Cells.Item(82, 11).Value = 0  # K82
$ws.Cells.Item(82, 12).Value = 0  # L82
$ws.Cells.Item(82, 13).ClearContents()  # M82
$ws.Cells.Item(82, 14).ClearContents()  # N82
# Row 85
$ws.Cells.Item(85, 8).Value = 0  # H85
$ws.Cells.Item(85, 9).Value = 0  # I85
$ws.Cells.Item(85, 10).Value = 0  # J85
$ws.Cells.Item(85, 11).Value = 0  # K85
$ws.Cells.Item(85, 12).Value = 0  # L85
$ws.Cells.Item(85, 13).ClearContents()  # M85
$ws.Cells.Item(85, 14).ClearContents()  # N85
# Row 109
$ws.Cells.Item(109, 8).Value = 42419.668  # H109
$ws.Cells.Item(109, 9).Value = 12259  # I109
$ws.Cells.Item(109, 10).Value = 57500  # J109
$ws.Cells.Item(109, 11).Value = 12259  # K109
$ws.Cells.Item(109, 12).Value = 57500  # L109
$ws.Cells.Item(109, 13).Value = -10872  # M109
$ws.Cells.Item(109, 14).Value = -60274  # N109
# Row 126
$ws.Cells.Item(126, 8).Value = 5092.12  # H126
$ws.Cells.Item(126, 9).Value = 4331.737  # I126
$ws.Cells.Item(126, 10).Value = 7500  # J126
$ws.Cells.Item(126, 11).Value = 12995.211  # K126
$ws.Cells.Item(126, 12).Value = 22500  # L126
$ws.Cells.Item(126, 13).Value = -10525.211  # M126
$ws.Cells.Item(126, 14).Value = -27440  # N126

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 74
$ws.Cells.Item(74, 8).Value = 11754.8  # H74
$ws.Cells.Item(74, 10).Value = 11754.8  # J74
$ws.Cells.Item(74, 12).Value = 11754.8  # L74
$ws.Cells.Item(74, 14).Value = -13626.8  # N74
# Row 77
$ws.Cells.Item(77, 8).Value = 11754.8  # H77
$ws.Cells.Item(77, 10).Value = 11754.8  # J77
$ws.Cells.Item(77, 12).Value = 35264.39999999999  # L77
$ws.Cells.Item(77, 14).Value = -44624.39999999999  # N77
# Row 113
$ws.Cells.Item(113, 8).Value = 3535.5144  # H113
$ws.Cells.Item(113, 9).Value = 4998.091  # I113
$ws.Cells.Item(113, 10).Value = 1060.3846  # J113
$ws.Cells.Item(113, 11).Value = 14994.273  # K113
$ws.Cells.Item(113, 12).Value = 3181.1538  # L113
$ws.Cells.Item(113, 13).Value = -12824.273  # M113
$ws.Cells.Item(113, 14).Value = -7521.1538  # N113
